$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.044999999999999
$ws.Range("C3").Value = -10.872
$ws.Range("D3").Value = -6.949
$ws.Range("C4").Value = -12.446
$ws.Range("D9").Value = -6.94
$ws.Range("B11").Value = 6.414
$ws.Range("B12").Value = 4.806
$ws.Range("C14").Value = -12.527
$ws.Range("B15").Value = 5.07
$ws.Range("D15").Value = -8.140000000000001
$ws.Range("D19").Value = -8.099
$ws.Range("D20").Value = -7.825
$ws.Range("D25").Value = -7.867
$ws.Range("C26").Value = -13.088
$ws.Range("B27").Value = 6.085000000000001
$ws.Range("D27").Value = -8.454000000000001
$ws.Range("B28").Value = 6.018
$ws.Range("D28").Value = -7.874000000000001
$ws.Range("D30").Value = -7.222999999999999
$ws.Range("B31").Value = 5.393000000000001
$ws.Range("C31").Value = -12.533
$ws.Range("B32").Value = 6.815
$ws.Range("D32").Value = -7.74
$ws.Range("C35").Value = -12.097
$ws.Range("B36").Value = 9.248000000000001
$ws.Range("C37").Value = -13.859
$ws.Range("B38").Value = 5.223999999999999
$ws.Range("C39").Value = -12.434
$ws.Range("C40").Value = -13.041
$ws.Range("D44").Value = -7.673999999999999
$ws.Range("C45").Value = -12.771
$ws.Range("B46").Value = 6.390000000000001
$ws.Range("D47").Value = -7.422
$ws.Range("C52").Value = -11.07
$ws.Range("B54").Value = 5.006
$ws.Range("B55").Value = 4.537
$ws.Range("B56").Value = 4.249000000000001
$ws.Range("C57").Value = -13.646
$ws.Range("D58").Value = -7.974000000000001
$ws.Range("D62").Value = -7.869
$ws.Range("B67").Value = 5.169
$ws.Range("B69").Value = 5.147
$ws.Range("B72").Value = 5.501
$ws.Range("B73").Value = 7.898999999999999
$ws.Range("D77").Value = -7.535000000000001
$ws.Range("D78").Value = -8.077
$ws.Range("C81").Value = -13.002
$ws.Range("B83").Value = 5.412000000000001
$ws.Range("C83").Value = -13.607
$ws.Range("D84").Value = -8.019
$ws.Range("B86").Value = 5.052000000000001
$ws.Range("D89").Value = -6.871
$ws.Range("B91").Value = 5.949
$ws.Range("D91").Value = -6.447
$ws.Range("D92").Value = -6.672
$ws.Range("B93").Value = 5.685999999999999
$ws.Range("D96").Value = -7.345000000000001
$ws.Range("B99").Value = 5.718
$ws.Range("C100").Value = -12.807
$ws.Range("C102").Value = -13.583
$ws.Range("D102").Value = -7.589
